# "Loan RBI, Variable Instalments"
#
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N -- this shifts the existing "Late" / "heading" (Date) /
# "Outstanding" columns from N/O/P to O/P/Q and widens the used range
# from A1:P15 to A1:Q15. The sheet also becomes the active tab (moving
# away from "Input").

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active/selected sheet (was "Input").
[void]$wsRepay.Activate()

# Insert a blank column before column N (14th column). Excel's default
# column-insert behaviour copies formatting from the column on the left
# (column M here), so give the new column that same width.
[void]$wsRepay.Columns.Item(14).Insert()
$wsRepay.Columns.Item(14).ColumnWidth = $wsRepay.Columns.Item(13).ColumnWidth

# Update the remembered selection on the sheet.
[void]$wsRepay.Range("R6").Select()
